# Daily attendance processing - 2025-10-11 21:40:35
# Swap the order of the first two comma-separated "Recorded By" entries
# in column G for the affected session rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,10,11,12,13,14,15,17,29,30,31,32,33,37,38,39,40,41,42,44,56,57,58,59,60,64,65,66,67,68,69,71,84,85,86,87,88,89,93,95,96,110,111,112,113,114,115,119,121,122,136,137,138,139,140,141,145,147,148)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value = [string]::Join(", ", $parts)
        }
    }
}
